# Generate Report for Handoff
#
# The localization-status report previously listed:
#   Overview / zh-cn / de-de  rows for:
#     row 2 -> 79d2a641-b11c-4d0e-b0b3-cee54dd48687.md  (Handed back: in sync with en-US)
#     row 3 -> 8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md  (Handed back: in sync with en-US)
#
# After re-running the handoff generation, the two files swapped rows and the
# 79d2a641... file now has a fresh handoff in progress ("Ready for handoff")
# while the 8c88d43a... file keeps its prior "Handed back" status:
#     row 2 -> 8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md  (Handed back: in sync with en-US)
#     row 3 -> 79d2a641-b11c-4d0e-b0b3-cee54dd48687.md  (Ready for handoff)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-37-20 08:37:04"

$ov.Range("A3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-38-20 08:38:03"

# ---------------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.50c414ecd8e910df6e3226df088baf65c337638d.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-20 08:37:00"
$zh.Range("F2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md"
$zh.Range("G2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.50c414ecd8e910df6e3226df088baf65c337638d.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-20 08:37:28"
$zh.Range("I2").Value = "Include"

$zh.Range("A3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.7c4116a52db88b6ac2402d952247cb9e1fa7bdf3.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-20 08:37:59"
$zh.Range("F3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.md"
$zh.Range("G3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.7c4116a52db88b6ac2402d952247cb9e1fa7bdf3.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-20 08:37:28"
$zh.Range("I3").Value = "Include"

# ---------------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.50c414ecd8e910df6e3226df088baf65c337638d.de-de.xlf"
$de.Range("E2").Value = "2016-03-20 08:37:04"
$de.Range("F2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md"
$de.Range("G2").Value = "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.50c414ecd8e910df6e3226df088baf65c337638d.de-de.xlf"
$de.Range("H2").Value = "2016-03-20 08:37:34"
$de.Range("I2").Value = "Include"

$de.Range("A3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.7c4116a52db88b6ac2402d952247cb9e1fa7bdf3.de-de.xlf"
$de.Range("E3").Value = "2016-03-20 08:38:03"
$de.Range("F3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.md"
$de.Range("G3").Value = "79d2a641-b11c-4d0e-b0b3-cee54dd48687.7c4116a52db88b6ac2402d952247cb9e1fa7bdf3.de-de.xlf"
$de.Range("H3").Value = "2016-03-20 08:37:34"
$de.Range("I3").Value = "Include"

Write-Host "Done updating localization-status report."
